# Applies the cryptos price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.095.22"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.889.66"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'245.09"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "'0.662"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'40.99"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").Value = "'0.346"
$ws.Range("E9").Value = "  +4.63%  "
$ws.Range("D10").Value = "'52.60"
$ws.Range("E10").Value = "  +12.01%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "'0.0991"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "2.164.31"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "'12.05"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "1.895.12"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'4.78"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "35.084.84"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'72.61"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "'239.39"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'12.35"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'4.79"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "'2.34"
$ws.Range("E26").Value = "  +21.98%  "
$ws.Range("D27").Value = "'170.00"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").Value = "'18.21"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").Value = "'4.09"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("B32").Value = "BinanceUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D32").Value = "'1.01"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0558"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'0.931"
$ws.Range("E34").Value = "  +13.48%  "
$ws.Range("D35").Value = "'4.04"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "'1.75"
$ws.Range("E36").Value = "  -4.97%  "
$ws.Range("D37").Value = "'2.01"
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "'0.0206"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("D41").Value = "'15.90"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("D42").Value = "'0.0618"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").Value = "'88.90"
$ws.Range("D44").Value = "1.333.76"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'47.67"
$ws.Range("E46").Value = "  +36.49%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "'2.77"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "'6.45"
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("D50").Value = "2.072.98"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").Value = "'11.51"
$ws.Range("E51").Value = "  -7.15%  "
